$p = $ppt.ActivePresentation

foreach ($i in 29..33) {
    $s = $p.Slides.Item($i)
    $s.SlideShowTransition.Hidden = $true
}
